$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1307.826
$ws.Range("I15").Value = 1307.826
$ws.Range("K15").Value = 3923.478
$ws.Range("M15").Value = -3754.478
$ws.Range("H76").Value = 3148.7
$ws.Range("I76").Value = 2749.8333
$ws.Range("J76").Value = 3747
$ws.Range("K76").Value = 2749.8333
$ws.Range("L76").Value = 3747
$ws.Range("M76").Value = -2434.8333
$ws.Range("N76").Value = -4377
$ws.Range("H79").Value = 3148.7
$ws.Range("I79").Value = 2749.8333
$ws.Range("J79").Value = 3747
$ws.Range("K79").Value = 2749.8333
$ws.Range("L79").Value = 3747
$ws.Range("M79").Value = -1657.8333
$ws.Range("N79").Value = -5931
$ws.Range("H98").Value = 305.93332
$ws.Range("I98").Value = 346.84616
$ws.Range("K98").Value = 346.84616
$ws.Range("M98").Value = 1151.15384
$ws.Range("H106").Value = 2358.75
$ws.Range("I106").Value = 2410
$ws.Range("K106").Value = 2410
$ws.Range("M106").Value = -1779
$ws.Range("H113").Value = 3725
$ws.Range("J113").Value = 3725
$ws.Range("L113").Value = 3725
$ws.Range("N113").Value = -10233
$ws.Range("H122").Value = 305.93332
$ws.Range("I122").Value = 346.84616
$ws.Range("K122").Value = 1040.53848
$ws.Range("M122").Value = 1409.46152

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22928.523
$ws.Range("I32").Value = 22286.086
$ws.Range("K32").Value = 22286.086
$ws.Range("M32").Value = -21999.086
$ws.Range("H45").Value = 3428.5715
$ws.Range("H102").Value = 2923
$ws.Range("I102").Value = 3047.6
$ws.Range("J102").Value = 2300
$ws.Range("K102").Value = 3047.6
$ws.Range("L102").Value = 2300
$ws.Range("M102").Value = -1425.6
$ws.Range("N102").Value = -5544
$ws.Range("H122").Value = 4655.485
$ws.Range("I122").Value = 3797.6206
$ws.Range("K122").Value = 11392.8618
$ws.Range("M122").Value = -8942.861800000001
$ws.Range("H132").Value = 16991.273
$ws.Range("I132").Value = 21433.178
$ws.Range("K132").Value = 64299.534
$ws.Range("M132").Value = -61769.534

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 619.4
$ws.Range("J64").Value = 611.75
$ws.Range("L64").Value = 611.75
$ws.Range("N64").Value = -1061.75
$ws.Range("H67").Value = 619.4
$ws.Range("J67").Value = 611.75
$ws.Range("L67").Value = 611.75
$ws.Range("N67").Value = -2171.75
$ws.Range("H86").Value = 3732.1052
$ws.Range("I86").Value = 2545.889
$ws.Range("J86").Value = 4799.7
$ws.Range("K86").Value = 2545.889
$ws.Range("L86").Value = 4799.7
$ws.Range("M86").Value = -1422.889
$ws.Range("N86").Value = -7045.7
$ws.Range("H89").Value = 3732.1052
$ws.Range("I89").Value = 2545.889
$ws.Range("J89").Value = 4799.7
$ws.Range("K89").Value = 12729.445
$ws.Range("L89").Value = 23998.5
$ws.Range("M89").Value = -7113.445
$ws.Range("N89").Value = -35230.5
$ws.Range("H107").Value = 2183.4375
$ws.Range("I107").Value = 1343.5714
$ws.Range("K107").Value = 1343.5714
$ws.Range("M107").Value = 576.4286
$ws.Range("H134").Value = 1848.425
$ws.Range("I134").Value = 1586.1471
$ws.Range("J134").Value = 3334.6667
$ws.Range("K134").Value = 4758.4413
$ws.Range("L134").Value = 10004.0001
$ws.Range("M134").Value = -2223.4413
$ws.Range("N134").Value = -15074.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 14327.77
$ws.Range("I99").Value = 22646.334
$ws.Range("J99").Value = 7197.5713
$ws.Range("K99").Value = 22646.334
$ws.Range("L99").Value = 7197.5713
$ws.Range("M99").Value = -21148.334
$ws.Range("N99").Value = -10193.5713
$ws.Range("H126").Value = 14327.77
$ws.Range("I126").Value = 22646.334
$ws.Range("J126").Value = 7197.5713
$ws.Range("K126").Value = 67939.00199999999
$ws.Range("L126").Value = 21592.7139
$ws.Range("M126").Value = -65469.00199999999
$ws.Range("N126").Value = -26532.7139
$ws.Range("H134").Value = 2461.9697
$ws.Range("I134").Value = 1374.381
$ws.Range("K134").Value = 4123.143
$ws.Range("M134").Value = -1588.143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 349.6
$ws.Range("I16").Value = 249.33333
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 747.99999
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -574.99999
$ws.Range("N16").Value = -1846
$ws.Range("H23").Value = 433.8
$ws.Range("I23").Value = 235
$ws.Range("J23").Value = 483.5
$ws.Range("K23").Value = 705
$ws.Range("L23").Value = 1450.5
$ws.Range("M23").Value = -470
$ws.Range("N23").Value = -1920.5
$ws.Range("H110").Value = 24999.75
$ws.Range("I110").Value = 24999
$ws.Range("K110").Value = 74997
$ws.Range("M110").Value = -70907
$ws.Range("H114").Value = 3285.625
$ws.Range("J114").Value = 3612.5715
$ws.Range("L114").Value = 10837.7145
$ws.Range("N114").Value = -17345.7145
$ws.Range("H121").Value = 1037.0769
$ws.Range("I121").Value = 223.66667
$ws.Range("J121").Value = 1281.1
$ws.Range("K121").Value = 671.00001
$ws.Range("L121").Value = 3843.3
$ws.Range("M121").Value = 638.99999
$ws.Range("N121").Value = -6463.299999999999
$ws.Range("H122").Value = 968.625
$ws.Range("I122").Value = 922.5
$ws.Range("J122").Value = 984
$ws.Range("K122").Value = 8302.5
$ws.Range("L122").Value = 8856
$ws.Range("M122").Value = -5852.5
$ws.Range("N122").Value = -13756
$ws.Range("H126").Value = 21666.334
$ws.Range("I126").Value = 14999
$ws.Range("K126").Value = 44997
$ws.Range("M126").Value = -40057

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8095.654
$ws.Range("I70").Value = 7888.4443
$ws.Range("K70").Value = 7888.4443
$ws.Range("M70").Value = -7618.4443
$ws.Range("H73").Value = 8095.654
$ws.Range("I73").Value = 7888.4443
$ws.Range("K73").Value = 7888.4443
$ws.Range("M73").Value = -6952.4443
$ws.Range("H102").Value = 16673447
$ws.Range("I102").Value = 23816430
$ws.Range("K102").Value = 23816430
$ws.Range("M102").Value = -23814808
$ws.Range("H126").Value = 3939.5293
$ws.Range("I126").Value = 2152.182
$ws.Range("K126").Value = 6456.545999999999
$ws.Range("M126").Value = -3986.545999999999
$ws.Range("H132").Value = 1656.9286
$ws.Range("I132").Value = 695.2
$ws.Range("K132").Value = 2085.6
$ws.Range("M132").Value = 444.3999999999996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2708.1875
$ws.Range("I16").Value = 1830.8
$ws.Range("K16").Value = 1830.8
$ws.Range("M16").Value = -1660.8
$ws.Range("H22").Value = 738.96295
$ws.Range("I22").Value = 698.4666999999999
$ws.Range("J22").Value = 789.5833
$ws.Range("K22").Value = 698.4666999999999
$ws.Range("L22").Value = 789.5833
$ws.Range("M22").Value = -403.4666999999999
$ws.Range("N22").Value = -1379.5833
$ws.Range("H27").Value = 738.96295
$ws.Range("I27").Value = 698.4666999999999
$ws.Range("J27").Value = 789.5833
$ws.Range("K27").Value = 698.4666999999999
$ws.Range("L27").Value = 789.5833
$ws.Range("M27").Value = -591.4666999999999
$ws.Range("N27").Value = -1003.5833
$ws.Range("H136").Value = 3649.1343
$ws.Range("I136").Value = 2671.06
$ws.Range("J136").Value = 6525.8237
$ws.Range("K136").Value = 8013.18
$ws.Range("L136").Value = 19577.4711
$ws.Range("M136").Value = -5463.18
$ws.Range("N136").Value = -24677.4711

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 99990
$ws.Range("J119").Value = 99990
$ws.Range("L119").Value = 99990
$ws.Range("N119").Value = -109666
$ws.Range("H126").Value = 2250.2222
$ws.Range("I126").Value = 1959.4
$ws.Range("K126").Value = 5878.200000000001
$ws.Range("M126").Value = -3408.200000000001
$ws.Range("H132").Value = 1698.8472
$ws.Range("I132").Value = 680.4737
$ws.Range("J132").Value = 2063.9246
$ws.Range("K132").Value = 2041.4211
$ws.Range("L132").Value = 6191.773799999999
$ws.Range("M132").Value = 488.5789
$ws.Range("N132").Value = -11251.7738
$ws.Range("H136").Value = 2967.6924
$ws.Range("I136").Value = 1964
$ws.Range("K136").Value = 5892
$ws.Range("M136").Value = -3342

Write-Output "done"